$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append a new data row (row 37) ---
# Continues the "Applications -> App 2 -> ... -> Applications -> App 20" pattern
# already used by the preceding rows 32-36 (Applications/App 2 depending on each
# of the other CI types in turn; this one adds the "Applications" dependency).
$descStart = $ws.Range("C36").Text
$descEnd   = $ws.Range("G36").Text

$rowVals = @("Applications", "App 2", $descStart, "Depends On", "Applications", "App 20", $descEnd)
$cols = @("A", "B", "C", "D", "E", "F", "G")

for ($i = 0; $i -lt $cols.Length; $i++) {
  $cell = $ws.Range($cols[$i] + "37")
  $cell.Value = $rowVals[$i]
  # Match formatting used by the rest of the table body (8pt font,
  # centered horizontally, top-aligned vertically).
  $cell.Font.Size = 8
  $cell.HorizontalAlignment = -4108   # xlCenter
  $cell.VerticalAlignment = -4160     # xlTop
}

# Columns C and G wrap their (longer) descriptive text, same as every other row.
$ws.Range("C37").WrapText = $true
$ws.Range("G37").WrapText = $true

# Row height matches the other (wrapped, 8pt) data rows.
$ws.Rows("37").RowHeight = $ws.Rows("36").RowHeight

# --- View/selection state ---
[void]$ws.Activate()
[void]$ws.Range("F30").Select()

# Best-effort: reflect the scrolled viewport (topLeftCell moved from A19 to A22)
# and the resized/maximized application window from the saved workbook view.
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$excel.Left = 38280
$excel.Top = 4080
$excel.Width = 29040
$excel.Height = 16440
